# This script applies the commit's data changes to the PDH (Setting Baju)
# mail-merge style document. Each field appears twice (once per "page"
# table cell), with identical old -> new text updates in both places.
#
# Word Find.Execute signature used below:
#   Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#           MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#           Format, ReplaceWith, Replace)
# Replace: 0 = wdReplaceNone, 1 = wdReplaceOne, 2 = wdReplaceAll

$d = $word.ActiveDocument

# --- Simple fields: unique text in the whole document (2 occurrences
#     each, identical substitution in both), safe to Replace All. ---

$d.Content.Find.Execute("Q34", $true, $false, $false, $false, $false, `
    $true, 1, $false, "20", 2) | Out-Null

$d.Content.Find.Execute("HADI SUDARMINTO", $true, $false, $false, $false, `
    $false, $true, 1, $false, "M. FARID SAMODRA", 2) | Out-Null

$d.Content.Find.Execute("DP 3 NAUTIKA / 30", $true, $false, $false, $false, `
    $false, $true, 1, $false, "DP 5 TEKNIK", 2) | Out-Null

# Numeric / short alpha fields - restrict to whole word matches so we
# don't clobber substrings elsewhere.
$d.Content.Find.Execute("41", $true, $true, $false, $false, $false, `
    $true, 1, $false, "42", 2) | Out-Null

$d.Content.Find.Execute("XL", $true, $true, $false, $false, $false, `
    $true, 1, $false, "M", 2) | Out-Null

$d.Content.Find.Execute("56", $true, $true, $false, $false, $false, `
    $true, 1, $false, "55", 2) | Out-Null

$d.Content.Find.Execute("40", $true, $true, $false, $false, $false, `
    $true, 1, $false, "46", 2) | Out-Null

$d.Content.Find.Execute("16", $true, $true, $false, $false, $false, `
    $true, 1, $false, "17", 2) | Out-Null

$d.Content.Find.Execute("70", $true, $true, $false, $false, $false, `
    $true, 1, $false, "66", 2) | Out-Null

$d.Content.Find.Execute("39", $true, $true, $false, $false, $false, `
    $true, 1, $false, "38", 2) | Out-Null

# --- "24" is ambiguous: it occurs twice per table cell (UB_2 and
#     UB_5), and they diverge (23 and 25 respectively), while the
#     other "26" occurrences (UB_4 / UB_6) must stay untouched.
#     Walk paragraph-by-paragraph, replacing the first "24" with "23"
#     and the following "24" (later in the same paragraph) with "25".

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i).Range
    if ($para.Text -like "*Uk*Baju*") {
        $paraEnd = $para.End

        $r1 = $d.Range($para.Start, $paraEnd)
        $found1 = $r1.Find.Execute("24", $true, $true, $false, $false, `
            $false, $true, 1, $false, "23", 1)

        if ($found1) {
            $r2 = $d.Range($r1.End, $paraEnd)
            $r2.Find.Execute("24", $true, $true, $false, $false, $false, `
                $true, 1, $false, "25", 1) | Out-Null
        }
    }
}
